# Daily update at 8 AM UTC
# Adds the new day's data row (45649 -> 2025-01-... ) and shifts the
# "most recent day" date-only formatting down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (62) loses its special "date only" formatting
# and reverts to the standard date+time format used by all other rows.
$ws.Range("A62").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 63.
$ws.Range("A63").Value2 = 45649
$ws.Range("B63").Value2 = 149
$ws.Range("C63").Value2 = 139
$ws.Range("D63").Value2 = 146

# The new last row gets the special "date only" formatting.
$ws.Range("A63").NumberFormat = "YYYY-MM-DD"
